# Regenerate save_data to use K (strikeouts) instead of Strike# for the G column.
# Column G header is "K"; values below are recalculated strikeout counts (s_vals).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 3
    3  = 3
    4  = 2
    5  = 2
    6  = 2
    7  = 4
    8  = 1
    9  = 1
    10 = 5
    11 = 2
    12 = 6
    13 = 4
    14 = 2
    15 = 3
    16 = 4
    17 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
